$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2023" column (K), mirroring the formatting
# already used by the preceding "2022" column (J).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1106
$ws.Range("K5").Value = 751.5
$ws.Range("K6").Value = 1245.0999999999999
